$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.706.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.794.15"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.14%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.31%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.791.18"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.17%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("E9").Value = "  +0.13%  "
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.460"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("E13").Value = "  +9.91%  "
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.430.55"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.843.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.798.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.65%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.42"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("E20").Value = "  +0.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.80"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "467.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.732"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("E24").Value = "  -8.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.20"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("E26").Value = "  +2.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.20"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.93%  "
$ws.Range("E29").Value = "  -0.10%  "
$ws.Range("E30").Value = "  -1.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.943.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.05%  "
$ws.Range("E32").Value = "  -3.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.69"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.26"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.760.63"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.36%  "
$ws.Range("E37").Value = "  +1.76%  "
$ws.Range("E38").Value = "  +6.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.139"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("E41").Value = "  -1.31%  "
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.317"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.49%  "
$ws.Range("E46").Value = "  -1.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "405.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -7.37%  "
$ws.Range("E49").Value = "  +0.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.92"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0358"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.09%  "
